# Add a new "range-all-cells" worksheet at the end of the workbook and
# populate it with a small header/data grid, matching the commit that
# adds the ability to extract a full table using the FROM_EXCEL
# transformer.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (so it ends up last:
# data, with separators, range-all-cells).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "range-all-cells"

# Header row with date columns (B2:D2).
$ws.Range("B2").Value = "date1"
$ws.Range("C2").Value = "date2"
$ws.Range("D2").Value = "date3"

# Data row for "Name 1".
$ws.Range("A3").Value = "Name 1"
$ws.Range("B3").Value = "line1"
$ws.Range("C3").Value = "line2"
$ws.Range("D3").Value = "line3"

# Data row for "Name 2".
$ws.Range("A4").Value = "Name 2"
$ws.Range("B4").Value = "line4"
$ws.Range("C4").Value = "line5"
$ws.Range("D4").Value = "line6"

# Match the saved selection on the new sheet.
$null = $ws.Range("E9").Select()
